$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# ---------------------------------------------------------------------------
# Tabelle1 ("EngineMerge" source sheet): insert two new flag columns
# (NETTO_MPS, MCRAD_MPS) right before the existing DispName/Title columns.
# Selecting AS:AT and inserting shifts the old AS/AT (DispName/Title) two
# columns to the right, becoming AU/AV, and the new AS/AT inherit the
# formatting of column AR (the existing VARIO_MPS flag column).
# ---------------------------------------------------------------------------
$ws1.Columns("AS:AT").Insert()

# Header row - claim shared-string slot 298 (NETTO_MPS) first so the
# shared-string table ends up in the same order as the source workbook.
$ws1.Range("AS1").Value = "NETTO_MPS"

# ---------------------------------------------------------------------------
# Tabelle2 ("EngineMerge" flattened/merge sheet): insert five new lookup
# columns (NETTO_MPS, NETTO_KTS, NETTO_ANI, MCRAD_MPS, MCRAD_KT) right
# before the old END_OF_COL/Title columns (which shift from EL/EM to EQ/ER).
# ---------------------------------------------------------------------------
$ws2.Columns("EL:EP").Insert()

$ws2.Range("EL1").Value = "NETTO_MPS"
$ws2.Range("EM1").Value = "NETTO_KTS"
$ws2.Range("EN1").Value = "NETTO_ANI"
$ws2.Range("EO1").Value = "MCRAD_MPS"
$ws2.Range("EP1").Value = "MCRAD_KT"

# Back to Tabelle1 - AT1 (MCRAD_MPS) now reuses the shared-string slot (301)
# that was claimed by Tabelle2's EO1 above.
$ws1.Range("AT1").Value = "MCRAD_MPS"

# ---------------------------------------------------------------------------
# Data rows 2-40 on Tabelle1: the two new columns mirror the existing
# VARIO_MPS flag column (AR) - 0 for every aircraft except the last row
# (40, "Generic Glider"), which is 1.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 40; $r++) {
    $flag = 0
    if ($r -eq 40) { $flag = 1 }
    $ws1.Cells.Item($r, 45).Value = $flag   # AS
    $ws1.Cells.Item($r, 46).Value = $flag   # AT
}

# ---------------------------------------------------------------------------
# Data rows 2-40 on Tabelle2: EL/EO pull the new Tabelle1 flags through a
# formula (like every other lookup column on this sheet); EM/EN/EP are
# filler columns (same "|" placeholder used throughout the sheet).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 40; $r++) {
    $ws2.Cells.Item($r, 142).Formula = "=Tabelle1!AS$r"   # EL
    $ws2.Cells.Item($r, 143).Value = "|"                  # EM
    $ws2.Cells.Item($r, 144).Value = "|"                  # EN
    $ws2.Cells.Item($r, 145).Formula = "=Tabelle1!AT$r"   # EO
    $ws2.Cells.Item($r, 146).Value = "|"                  # EP
}

# ---------------------------------------------------------------------------
# Workbook-level AutoFilter defined name: the filter range grows by two
# columns (AS -> AU) to account for the two newly inserted Tabelle1 columns.
# ---------------------------------------------------------------------------
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Tabelle1!`$A`$1:`$AU`$1"

# ---------------------------------------------------------------------------
# Restore the selections recorded in the source workbook after the edit.
# ---------------------------------------------------------------------------
$ws1.Range("AU12").Select()
$ws2.Range("AR50").Select()
